# Auto-generated edit script applying scheduled market-data refresh to Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 855.2222
$ws.Range("I129").Value = 657.4
$ws.Range("J129").Value = 869.98505
$ws.Range("K129").Value = 1972.2
$ws.Range("L129").Value = 2609.95515
$ws.Range("M129").Value = 3027.8
$ws.Range("N129").Value = -12609.95515
$ws.Range("H135").Value = 528.4666999999999
$ws.Range("I135").Value = 531.2143
$ws.Range("K135").Value = 4780.928699999999
$ws.Range("M135").Value = -2245.928699999999
$ws.Range("H138").Value = 1727.7819
$ws.Range("I138").Value = 1248.7587
$ws.Range("J138").Value = 2262.077
$ws.Range("K138").Value = 3746.2761
$ws.Range("L138").Value = 6786.231000000001
$ws.Range("M138").Value = 1393.7239
$ws.Range("N138").Value = -17066.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3496.9524
$ws.Range("I32").Value = 2127.7114
$ws.Range("J32").Value = 9969.727999999999
$ws.Range("K32").Value = 2127.7114
$ws.Range("L32").Value = 9969.727999999999
$ws.Range("M32").Value = -1840.7114
$ws.Range("N32").Value = -10543.728
$ws.Range("H61").Value = 3565.5
$ws.Range("I61").Value = 2598.4211
$ws.Range("K61").Value = 2598.4211
$ws.Range("M61").Value = -2386.4211
$ws.Range("H74").Value = 1837.238
$ws.Range("I74").Value = 1361.6875
$ws.Range("J74").Value = 3359
$ws.Range("K74").Value = 1361.6875
$ws.Range("L74").Value = 3359
$ws.Range("M74").Value = -487.6875
$ws.Range("N74").Value = -5107
$ws.Range("H77").Value = 1837.238
$ws.Range("I77").Value = 1361.6875
$ws.Range("J77").Value = 3359
$ws.Range("K77").Value = 6808.4375
$ws.Range("L77").Value = 16795
$ws.Range("M77").Value = -2440.4375
$ws.Range("N77").Value = -25531
$ws.Range("H110").Value = 1364.6207
$ws.Range("I110").Value = 1086.3846
$ws.Range("J110").Value = 3776
$ws.Range("K110").Value = 1086.3846
$ws.Range("L110").Value = 3776
$ws.Range("M110").Value = 958.6153999999999
$ws.Range("N110").Value = -7866
$ws.Range("H122").Value = 2199.5
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H136").Value = 3565.5
$ws.Range("I136").Value = 2598.4211
$ws.Range("K136").Value = 7795.263300000001
$ws.Range("M136").Value = -5245.263300000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 80011
$ws.Range("J18").Value = 80011
$ws.Range("L18").Value = 80011
$ws.Range("N18").Value = -81069
$ws.Range("H134").Value = 3159.55
$ws.Range("I134").Value = 3159.55
$ws.Range("K134").Value = 9478.650000000001
$ws.Range("M134").Value = -6943.650000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 450
$ws.Range("I2").Value = 450
$ws.Range("K2").Value = 450
$ws.Range("M2").Value = -337
$ws.Range("H7").Value = 375
$ws.Range("J7").Value = 1000
$ws.Range("L7").Value = 1000
$ws.Range("N7").Value = -1226
$ws.Range("H22").Value = 1207.625
$ws.Range("J22").Value = 1323
$ws.Range("L22").Value = 1323
$ws.Range("N22").Value = -2023
$ws.Range("H31").Value = 2020
$ws.Range("I31").Value = 1585.65
$ws.Range("K31").Value = 1585.65
$ws.Range("M31").Value = -1290.65
$ws.Range("H34").Value = 2020
$ws.Range("I34").Value = 1585.65
$ws.Range("K34").Value = 1585.65
$ws.Range("M34").Value = -1383.65
$ws.Range("H107").Value = 670.34784
$ws.Range("I107").Value = 681
$ws.Range("K107").Value = 681
$ws.Range("M107").Value = 1239
$ws.Range("H122").Value = 6253.5
$ws.Range("J122").Value = 5338
$ws.Range("L122").Value = 16014
$ws.Range("N122").Value = -20914
$ws.Range("H132").Value = 2076
$ws.Range("I132").Value = 1188
$ws.Range("K132").Value = 3564
$ws.Range("M132").Value = -1034
$ws.Range("H134").Value = 1180.4
$ws.Range("I134").Value = 1180.4
$ws.Range("K134").Value = 3541.2
$ws.Range("M134").Value = -1006.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 866.44446
$ws.Range("J5").Value = 914
$ws.Range("L5").Value = 2742
$ws.Range("N5").Value = -2966
$ws.Range("H105").Value = 2949
$ws.Range("J105").Value = 3063.5
$ws.Range("L105").Value = 9190.5
$ws.Range("N105").Value = -14432.5
$ws.Range("H122").Value = 1034
$ws.Range("I122").Value = 820.75
$ws.Range("J122").Value = 1105.0834
$ws.Range("K122").Value = 7386.75
$ws.Range("L122").Value = 9945.750599999999
$ws.Range("M122").Value = -4936.75
$ws.Range("N122").Value = -14845.7506
$ws.Range("H131").Value = 769.9
$ws.Range("J131").Value = 788.40216
$ws.Range("L131").Value = 2365.20648
$ws.Range("N131").Value = -12445.20648
$ws.Range("H135").Value = 866.44446
$ws.Range("J135").Value = 914
$ws.Range("L135").Value = 8226
$ws.Range("N135").Value = -13296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1483.1177
$ws.Range("I97").Value = 1450.25
$ws.Range("K97").Value = 1450.25
$ws.Range("M97").Value = -954.25
$ws.Range("H98").Value = 18206
$ws.Range("J98").Value = 18206
$ws.Range("L98").Value = 18206
$ws.Range("N98").Value = -24196
$ws.Range("H102").Value = 3089.9167
$ws.Range("I102").Value = 3531.111
$ws.Range("J102").Value = 1766.3334
$ws.Range("K102").Value = 3531.111
$ws.Range("L102").Value = 1766.3334
$ws.Range("M102").Value = -1909.111
$ws.Range("N102").Value = -5010.3334
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4239.5
$ws.Range("I7").Value = 2611.75
$ws.Range("J7").Value = 6409.8335
$ws.Range("K7").Value = 2611.75
$ws.Range("L7").Value = 6409.8335
$ws.Range("M7").Value = -2499.75
$ws.Range("N7").Value = -6633.8335
$ws.Range("H16").Value = 4192.154
$ws.Range("I16").Value = 6580
$ws.Range("K16").Value = 6580
$ws.Range("M16").Value = -6410
$ws.Range("H43").Value = 13205.6
$ws.Range("J43").Value = 13205.6
$ws.Range("L43").Value = 13205.6
$ws.Range("N43").Value = -13591.6
$ws.Range("H55").Value = 250.3125
$ws.Range("I55").Value = 205.46153
$ws.Range("J55").Value = 444.66666
$ws.Range("K55").Value = 205.46153
$ws.Range("L55").Value = 444.66666
$ws.Range("M55").Value = -32.46153000000001
$ws.Range("N55").Value = -790.66666
$ws.Range("H61").Value = 2424
$ws.Range("I61").Value = 1565.3334
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 1565.3334
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -1363.3334
$ws.Range("N61").Value = -5404
$ws.Range("H82").Value = 1249.6666
$ws.Range("I82").Value = 1339.8
$ws.Range("J82").Value = 799
$ws.Range("K82").Value = 1339.8
$ws.Range("L82").Value = 799
$ws.Range("M82").Value = -978.8
$ws.Range("N82").Value = -1521
$ws.Range("H85").Value = 1249.6666
$ws.Range("I85").Value = 1339.8
$ws.Range("J85").Value = 799
$ws.Range("K85").Value = 1339.8
$ws.Range("L85").Value = 799
$ws.Range("M85").Value = -91.79999999999995
$ws.Range("N85").Value = -3295
$ws.Range("H93").Value = 916.6667
$ws.Range("I93").Value = 975
$ws.Range("J93").Value = 800
$ws.Range("K93").Value = 975
$ws.Range("L93").Value = 800
$ws.Range("M93").Value = 273
$ws.Range("N93").Value = -3296
$ws.Range("H113").Value = 2424
$ws.Range("I113").Value = 1565.3334
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 1565.3334
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 604.6666
$ws.Range("N113").Value = -9340
$ws.Range("H122").Value = 7831.273
$ws.Range("I122").Value = 6879.125
$ws.Range("J122").Value = 10370.333
$ws.Range("K122").Value = 20637.375
$ws.Range("L122").Value = 31110.999
$ws.Range("M122").Value = -18187.375
$ws.Range("N122").Value = -36010.999
$ws.Range("H126").Value = 4239.5
$ws.Range("I126").Value = 2611.75
$ws.Range("J126").Value = 6409.8335
$ws.Range("K126").Value = 7835.25
$ws.Range("L126").Value = 19229.5005
$ws.Range("M126").Value = -5365.25
$ws.Range("N126").Value = -24169.5005
$ws.Range("H132").Value = 2875.84
$ws.Range("I132").Value = 2056.1875
$ws.Range("K132").Value = 6168.5625
$ws.Range("M132").Value = -3638.5625
$ws.Range("H136").Value = 3382.9
$ws.Range("J136").Value = 3748.2222
$ws.Range("L136").Value = 11244.6666
$ws.Range("N136").Value = -16344.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13999
$ws.Range("J41").Value = 13999
$ws.Range("L41").Value = 13999
$ws.Range("N41").Value = -14779
$ws.Range("H45").Value = 17784.5
$ws.Range("I45").Value = 20569
$ws.Range("J45").Value = 15000
$ws.Range("K45").Value = 20569
$ws.Range("L45").Value = 15000
$ws.Range("M45").Value = -20078
$ws.Range("N45").Value = -15982
$ws.Range("H112").Value = 14500
$ws.Range("J112").Value = 14500
$ws.Range("L112").Value = 14500
$ws.Range("N112").Value = -17454
$ws.Range("H136").Value = 24156660
$ws.Range("I136").Value = 26456914
$ws.Range("K136").Value = 79370742
$ws.Range("M136").Value = -79368192
